$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("A1").Value = 2504.3641967919812
$ws.Range("B1").Value = 1675.4831224905172
$ws.Range("C1").Value = 1659.1742236099808
$ws.Range("A2").Value = 2263.2649178670267
$ws.Range("B2").Value = 1537.5656129537554
$ws.Range("C2").Value = 1393.3171188274771
$ws.Range("A3").Value = 2560.3680356546965
$ws.Range("B3").Value = 1740.678815942545
$ws.Range("C3").Value = 1579.3443776225236
$ws.Range("A4").Value = 2495.9981556247999
$ws.Range("B4").Value = 1913.0014038613795
$ws.Range("C4").Value = 1950.7277354949008
$ws.Range("A5").Value = 2528.8694636736568
$ws.Range("B5").Value = 1742.3270368293149
$ws.Range("C5").Value = 1784.7568909152653
$ws.Range("A6").Value = 2480.0523703440967
$ws.Range("B6").Value = 1850.3730377169611
$ws.Range("C6").Value = 1931.0848430119481
$ws.Range("A7").Value = 2397.2088800503925
$ws.Range("B7").Value = 1848.4300944536446
$ws.Range("C7").Value = 1675.668040406772
$ws.Range("A8").Value = 2467.9969983108795
$ws.Range("B8").Value = 1930.0320208708331
$ws.Range("C8").Value = 1785.2136052421129
$ws.Range("A9").Value = 2643.8787318127324
$ws.Range("B9").Value = 1963.5067209681306
$ws.Range("C9").Value = 1680.4583613843849
$ws.Range("A10").Value = 2394.9101549394882
$ws.Range("B10").Value = 1520.9883557933742
$ws.Range("C10").Value = 1428.5183673137938
$ws.Range("A11").Value = 2163.0420355061924
$ws.Range("B11").Value = 1575.4241743753942
$ws.Range("C11").Value = 1407.2608779610723
$ws.Range("A12").Value = 2770.6265500970344
$ws.Range("B12").Value = 2207.4257872533763
$ws.Range("C12").Value = 1916.0906809257067
$ws.Range("A13").Value = 2528.2387480558714
$ws.Range("B13").Value = 1947.5474619142117
$ws.Range("C13").Value = 1749.6306059708488
$ws.Range("A14").Value = 2609.7083426384261
$ws.Range("B14").Value = 2019.1228183797468
$ws.Range("C14").Value = 1756.8907062203525
$ws.Range("A15").Value = 2503.2255475334637
$ws.Range("B15").Value = 2051.809549640262
$ws.Range("C15").Value = 1832.9655400561949
$ws.Range("A16").Value = 2592.375598703285
$ws.Range("B16").Value = 1804.6455680313004
$ws.Range("C16").Value = 1530.1185850923541
$ws.Range("A17").Value = 2356.8601578997373
$ws.Range("B17").Value = 1814.1810148436214
$ws.Range("C17").Value = 1709.7047242962224
$ws.Range("A18").Value = 2615.4902880995196
$ws.Range("B18").Value = 2163.225143496114
$ws.Range("C18").Value = 2077.8195710672057
$ws.Range("A19").Value = 2011.2280643505962
$ws.Range("B19").Value = 2033.9321272375284
$ws.Range("C19").Value = 1981.0696108332579
$ws.Range("A20").Value = 2576.2845510916013
$ws.Range("B20").Value = 1981.3683885098346
$ws.Range("C20").Value = 1870.0663422816451
$ws.Range("A21").Value = 2713.0801993379068
$ws.Range("B21").Value = 2020.8570131593781
$ws.Range("C21").Value = 1980.6742669114219
$ws.Range("A22").Value = 2569.5114392622113
$ws.Range("B22").Value = 1982.4395818578371
$ws.Range("C22").Value = 1731.9556707015938
